# CryCompanywiseStockReport_1.xlsx — stock-take correction
#
# A number of line items had their counted Quantity (column F) reduced
# (shrinkage / re-count adjustment). For each corrected row the Stock
# Value (column G) must be re-derived as Rate (D) * Quantity (F), and
# every "Sub Total:" row must be re-derived as the sum of the Stock
# Value cells belonging to its company block; the overall "Sub Total:"
# row (sum of every company subtotal) and the final "Grand Total:" row
# must follow the same cascade.
#
# Columns: A=S.No  B=Code  C=Description  D=Rate  E=MRP  F=Qty  G=Value

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colA = 1; $colB = 2; $colD = 4; $colE = 5; $colF = 6; $colG = 7

# ---------------------------------------------------------------------
# 1) New quantities (row -> new Qty) for every corrected line item.
# ---------------------------------------------------------------------
$F_NEW = @{
    10 = 3;    32 = 31;   38 = 56;   87 = 33;   91 = 76;   96 = 56;
    100 = 272; 105 = 208; 110 = 103; 121 = 82;  132 = 10;  154 = 12;
    155 = 80;  157 = 127; 270 = 69;  282 = 34;  287 = 34;  296 = 16;
    303 = 80;  304 = 165; 305 = 1;   323 = 138; 334 = 150; 344 = 107;
    346 = 41;  361 = 14;  415 = 28;  433 = 48;  434 = 32;  439 = 85;
    466 = 593; 468 = 599; 469 = 277; 470 = 202; 473 = 372; 474 = 410;
    477 = 1000;480 = 433; 481 = 344; 483 = 610; 527 = 621; 533 = 136;
    544 = 77;  547 = 28;  568 = 20;  569 = 48;  570 = 174; 572 = 12;
    580 = 18;  582 = 16;  588 = 30;  590 = 42;  597 = 110; 600 = 60;
    605 = 80;  608 = 66;  609 = 62;  610 = 9;   615 = 0;   633 = 170;
    637 = 169; 697 = 498; 698 = 145; 700 = 97;  701 = 99;  702 = 497;
    711 = 37;  712 = 7;   715 = 992; 716 = 99;  718 = 98;  730 = 16;
    740 = 24;  742 = 71;  744 = 80;  749 = 3;   790 = 740; 794 = 62
}

# Rows 304/305: the stock code (B) together with its Rate (D) and MRP
# (E) were swapped between these two adjoining lines (same item,
# re-attributed between two receipts) in addition to the quantity
# change above.
$B_SWAP = @{ 304 = 61610; 305 = 57077 }
$D_SWAP = @{ 304 = 102.71; 305 = 93.08 }
$E_SWAP = @{ 304 = 122.71; 305 = 111.2 }

foreach ($row in $B_SWAP.Keys) {
    $ws.Cells.Item($row, $colB).Value = $B_SWAP[$row]
    $ws.Cells.Item($row, $colD).Value = $D_SWAP[$row]
    $ws.Cells.Item($row, $colE).Value = $E_SWAP[$row]
}

# ---------------------------------------------------------------------
# 2) Apply the new quantities and recompute Stock Value = Rate * Qty.
# ---------------------------------------------------------------------
foreach ($row in $F_NEW.Keys) {
    $qty = $F_NEW[$row]
    $ws.Cells.Item($row, $colF).Value = $qty
    $rate = $ws.Cells.Item($row, $colD).Value2
    $ws.Cells.Item($row, $colG).Value = $rate * $qty
}

# ---------------------------------------------------------------------
# 3) Recompute every "Sub Total:" row as the sum of the Stock Value
#    column over the data rows of its company block (the rows between
#    it and the previous "Sub Total:"/company-header boundary). The
#    final "Sub Total:" row (just above "Grand Total:") instead sums
#    every other "Sub Total:" row, and "Grand Total:" mirrors it.
# ---------------------------------------------------------------------
$lastRow = $ws.UsedRange.Rows.Count

$subtotalRows = New-Object System.Collections.ArrayList
for ($r = 1; $r -le $lastRow; $r++) {
    $label = $ws.Cells.Item($r, $colA).Value2
    if ($label -is [string] -and $label -eq "Sub Total:") {
        [void]$subtotalRows.Add($r)
    }
}

$prevBoundary = 0
foreach ($sr in $subtotalRows) {
    # Walk upward from the row above the subtotal until we hit a
    # company-header row (text in col A, blank col B) or the previous
    # block boundary.
    $start = $sr - 1
    while ($start -gt $prevBoundary) {
        $aVal = $ws.Cells.Item($start, $colA).Value2
        $bVal = $ws.Cells.Item($start, $colB).Value2
        $isHeader = ($aVal -is [string]) -and ($aVal -ne "Sub Total:") -and ($bVal -eq "")
        if ($isHeader) { break }
        $start--
    }
    $dataStart = $start + 1
    $dataEnd = $sr - 1

    if ($dataStart -le $dataEnd) {
        $sum = 0.0
        $sawSubtotal = $false
        for ($r = $dataStart; $r -le $dataEnd; $r++) {
            $aVal = $ws.Cells.Item($r, $colA).Value2
            if ($aVal -is [string] -and $aVal -eq "Sub Total:") { $sawSubtotal = $true }
        }
        if ($sawSubtotal) {
            # This subtotal's "block" is actually the set of previously
            # computed company subtotals (the final rollup row) — sum
            # column B over those rows instead of column G.
            for ($r = $dataStart; $r -le $dataEnd; $r++) {
                $sum += $ws.Cells.Item($r, $colB).Value2
            }
        } else {
            for ($r = $dataStart; $r -le $dataEnd; $r++) {
                $sum += $ws.Cells.Item($r, $colG).Value2
            }
        }
        $ws.Cells.Item($sr, $colB).Value = $sum
    } else {
        # No rows directly above (degenerate block): this is the
        # rollup "Sub Total:" row — sum every other "Sub Total:" row.
        $sum = 0.0
        foreach ($other in $subtotalRows) {
            if ($other -ne $sr) { $sum += $ws.Cells.Item($other, $colB).Value2 }
        }
        $ws.Cells.Item($sr, $colB).Value = $sum
    }

    $prevBoundary = $sr
}

# ---------------------------------------------------------------------
# 4) "Grand Total:" row mirrors the rollup "Sub Total:" row directly
#    above it.
# ---------------------------------------------------------------------
for ($r = 1; $r -le $lastRow; $r++) {
    $label = $ws.Cells.Item($r, $colA).Value2
    if ($label -is [string] -and $label -eq "Grand Total:") {
        $ws.Cells.Item($r, $colB).Value = $ws.Cells.Item($r - 1, $colB).Value2
    }
}
